# Apply the "nn results 2020" update:
#  - Row 14: I14 shared-formula range grows to I14:I15 (formula/value unchanged),
#            J14 switches from "41,3k" to "28k", L14 switches from 20 to 50.
#  - New row 15: another "Conv, (128, 128), (64, 128), (1,128)" run (epoch 1).
#  - New row 16: another such run (epoch 2), noted "back to full dataset".
#  - New row 17: a "Conv, (128, 130), (64, 128), (1,128)" run stub (no metrics yet).
#  - Active-cell selection moves from G23 to C20.
#
# New shared strings must be introduced in the same order they first appear
# in the finished workbook (28k, Conv(128,130).., back to full dataset, 57,4k)
# so they land at shared-string indices 35-38, matching the target file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 14 edits (introduces new string "28k") --------------------------
$ws.Range("J14").Value = "28k"
$ws.Range("L14").Value = 50

# --- Row 15 (new, reuses existing strings only) --------------------------
$ws.Range("A15").Value = 384
$ws.Range("B15").Value = 128
$ws.Range("C15").Value = "Conv, (128, 128), (64, 128), (1,128)"
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0.31
$ws.Range("F15").Value = 0.87
$ws.Range("G15").Value = 0.7
$ws.Range("H15").Value = 0.53
$ws.Range("I15").Formula = "=2*(G15*H15)/(G15+H15)"
$ws.Range("J15").Value = "41,3k"
$ws.Range("L15").Value = 45
$ws.Range("M15").Value = "very reduced dataset"

# --- Row 16 numeric/reused-string cells first -----------------------------
$ws.Range("A16").Value = 384
$ws.Range("B16").Value = 128
$ws.Range("C16").Value = "Conv, (128, 128), (64, 128), (1,128)"
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0.3
$ws.Range("F16").Value = 0.866
$ws.Range("G16").Value = 0.656
$ws.Range("H16").Value = 0.646
$ws.Range("I16").Formula = "=2*(G16*H16)/(G16+H16)"
$ws.Range("J16").Value = "41,3k"
$ws.Range("L16").Value = 20

# --- Row 17 structure cells (introduces new string "Conv, (128, 130)...")
$ws.Range("A17").Value = 384
$ws.Range("B17").Value = 128
$ws.Range("C17").Value = "Conv, (128, 130), (64, 128), (1,128)"
$ws.Range("D17").Value = 2
$ws.Range("L17").Value = 20

# --- Back to row 16's note (introduces new string "back to full dataset")
$ws.Range("M16").Value = "back to full dataset"

# --- Finally row 17's weight note (introduces new string "57,4k") --------
$ws.Range("J17").Value = "57,4k"

# --- Selection moves to C20 ---------------------------------------------
[void]$ws.Activate()
$ws.Range("C20").Select() | Out-Null

Write-Host "Edits applied."
